$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artifacts")

# Insert a new column before column E ("Architecture"). This shifts the old
# E:J ("Architecture".."Remarks") one to the right, to F:K, carrying along
# their existing (already-correct) values/types untouched.
$ws.Columns.Item(5).Insert()

# New column header: "Inapplicable CPE URIs"
$ws.Range("E1").Value = "Inapplicable CPE URIs"

# The existing "CPE URIs" column (D) used to hold one concatenated string
# for every row; now each row gets just the CPE that actually applies.
$ws.Range("D2").Value = "cpe:/o:microsoft:windows_11_25h2"
$ws.Range("D3").Value = "cpe:/o:microsoft:windows_11_24h2"
$ws.Range("D4").Value = "cpe:/o:microsoft:windows_11_23h2"

# Fill the new "Inapplicable CPE URIs" column.
$ws.Range("E2").Value = "cpe:/o:microsoft:windows_11"
$ws.Range("E3").Value = "cpe:/o:microsoft:windows_11"
$ws.Range("E4").Value = "cpe:/o:microsoft:windows_11"

# Fill the (previously empty) "Remarks" column, now K.
$remark = "cpe:/o:microsoft:windows_11 is considered a outdated for this version of windows."
$ws.Range("K2").Value = $remark
$ws.Range("K3").Value = $remark
$ws.Range("K4").Value = $remark

# Column widths (best-fit-ish, matching the refreshed layout).
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(2).ColumnWidth = 11.83
$ws.Columns.Item(3).ColumnWidth = 13.83
$ws.Columns.Item(4).ColumnWidth = 27.67
$ws.Columns.Item(5).ColumnWidth = 22.67
$ws.Columns.Item(6).ColumnWidth = 12.33
$ws.Columns.Item(7).ColumnWidth = 13.67
$ws.Columns.Item(8).ColumnWidth = 13
$ws.Columns.Item(9).ColumnWidth = 7.5
$ws.Columns.Item(10).ColumnWidth = 17.17
$ws.Columns.Item(11).ColumnWidth = 62.33

# Selection / view tidy-up to match the refreshed sheet.
$ws.Range("E2").Select()

# Keep the autofilter / defined name in sync with the new last column (K).
$ws.Range("A1:K65001").AutoFilter()

$wb.Windows.Item(1).WindowState = -4143
